# Updated cryptos list (Price + Volume(1h) columns, plus a Uniswap/Chainlink row swap)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''67.858.37'
$ws.Range("E2").Value = '  +0.80%  '
$ws.Range("D3").Value = '''2.493.65'
$ws.Range("E3").Value = '  -0.34%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '''587.03'
$ws.Range("E5").Value = '  +0.02%  '
$ws.Range("D6").Value = '''177.05'
$ws.Range("E6").Value = '  +2.42%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  +0.14%  '
$ws.Range("E9").Value = '  +3.32%  '
$ws.Range("E10").Value = '  -0.06%  '
$ws.Range("E11").Value = '  +1.60%  '
$ws.Range("E12").Value = '  -0.23%  '
$ws.Range("D13").Value = '''2.948.94'
$ws.Range("E13").Value = '  +1.00%  '
$ws.Range("D14").Value = '''25.66'
$ws.Range("E14").Value = '  +0.31%  '
$ws.Range("D15").Value = '''67.725.79'
$ws.Range("E15").Value = '  +1.18%  '
$ws.Range("D16").Value = '''0.0000172'
$ws.Range("E16").Value = '  +0.17%  '
$ws.Range("D17").Value = '''2.492.67'
$ws.Range("E17").Value = '  +1.88%  '
$ws.Range("B18").Value = 'Chainlink'
$ws.Range("C18").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D18").Value = '''10.97'
$ws.Range("E18").Value = '  -0.91%  '
$ws.Range("B19").Value = 'Uniswap'
$ws.Range("C19").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D19").Value = '''7.50'
$ws.Range("E19").Value = '  +0.47%  '
$ws.Range("D20").Value = '''350.83'
$ws.Range("E20").Value = '  -0.17%  '
$ws.Range("E21").Value = '  +2.28%  '
$ws.Range("E22").Value = '  -0.06%  '
$ws.Range("D23").Value = '''70.87'
$ws.Range("E23").Value = '  +3.21%  '
$ws.Range("D24").Value = '''4.30'
$ws.Range("E24").Value = '  +0.96%  '
$ws.Range("D25").Value = '''1.75'
$ws.Range("E25").Value = '  -3.28%  '
$ws.Range("D26").Value = '''9.13'
$ws.Range("E26").Value = '  -1.81%  '
$ws.Range("D27").Value = '''2.619.58'
$ws.Range("E27").Value = '  +0.02%  '
$ws.Range("D28").Value = '''0.998'
$ws.Range("E28").Value = '  -0.01%  '
$ws.Range("D29").Value = '''0.0₃0905'
$ws.Range("E29").Value = '  -0.56%  '
$ws.Range("D30").Value = '''504.45'
$ws.Range("E30").Value = '  -1.79%  '
$ws.Range("D31").Value = '''7.84'
$ws.Range("E31").Value = '  +0.22%  '
$ws.Range("E32").Value = '  +1.57%  '
$ws.Range("E33").Value = '  -0.20%  '
$ws.Range("E35").Value = '  +3.00%  '
$ws.Range("D36").Value = '''162.26'
$ws.Range("E36").Value = '  +1.54%  '
$ws.Range("D37").Value = '''18.66'
$ws.Range("E37").Value = '  -0.27%  '
$ws.Range("E38").Value = '  +0.14%  '
$ws.Range("D39").Value = '''1.34'
$ws.Range("E39").Value = '  -0.42%  '
$ws.Range("E40").Value = '  -0.05%  '
$ws.Range("E41").Value = '  +3.02%  '
$ws.Range("E42").Value = '  +0.15%  '
$ws.Range("E43").Value = '  -0.53%  '
$ws.Range("E44").Value = '  +1.14%  '
$ws.Range("D45").Value = '''144.65'
$ws.Range("E45").Value = '  +0.81%  '
$ws.Range("E46").Value = '  +1.44%  '
$ws.Range("D47").Value = '''0.514'
$ws.Range("E47").Value = '  -0.98%  '
$ws.Range("E48").Value = '  +0.60%  '
$ws.Range("D49").Value = '''0.0743'
$ws.Range("E49").Value = '  +1.65%  '
$ws.Range("D50").Value = '''1.58'
$ws.Range("E50").Value = '  +0.16%  '
$ws.Range("E51").Value = '  -0.14%  '
